$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Gennaro Bullo"
$ws.Range("B19").Value = "Raffaele Prezzi  | Hellas Lazio"
$ws.Range("C19").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D19").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("E19").Value = "Filippo Benetti | I Magnifici"
$ws.Range("F19").Value = "Mattia Bertolini | QUEI STRAZI"
